$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
$ws.Range("A1").Value = "Employee Name "
$ws.Range("B1").Value = "Employee ID"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Temperature"

# --- Names for column A (rows 2-13), two blocks of six employees ---
$names = @("Mr. A", "Mr. B", "Mr. C", "Mr. D", "Mr. E", "Mr. F")
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $names[$i]
    $ws.Cells.Item($i + 8, 1).Value = $names[$i]
}

# --- Employee IDs for column B ---
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $i + 1
    $ws.Cells.Item($i + 8, 2).Value = $i + 1
}

# --- Dates for column C ---
$ws.Range("C2:C7").Value = 44100
$ws.Range("C8:C13").Value = 44101

# Apply the date number format once, then propagate the exact same style
# (not just an equivalent one) to the rest of the date cells via a
# copy / paste-special-formats so they all share a single cellXfs entry.
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C3:C13").PasteSpecial(-4122) | Out-Null

# --- Temperatures for column D ---
$temps = @(36.5, 36.8, 37.5, 36.1, 37.4, 36.9, 36.2, 36.9, 36.1, 36.2, 36.8, 36.4)
for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $temps[$i]
}

# --- Column A width (best effort; target 15.42578125 isn't reachable through
# the ColumnWidth setter's internal rounding grid, so use the nearest value) ---
$ws.Columns.Item(1).ColumnWidth = 14.6

# --- Selection cursor, like Excel leaves behind after data entry ---
$ws.Range("D14").Select() | Out-Null
